# 2022 April Sprint Test cases - add "RPA Test Cases (alpha)" sheet,
# tweak a couple of values, and append version/ratio rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Content fix already on "RPA Test Cases": row 9 col F text changed.
# ---------------------------------------------------------------------------
$ws1.Range("F9").Value = "Success  queue and CCR "

# ---------------------------------------------------------------------------
# 2) Duplicate the sheet to create "RPA Test Cases (alpha)" right after it.
# ---------------------------------------------------------------------------
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "RPA Test Cases (alpha)"

# ---------------------------------------------------------------------------
# 3) On the new sheet, clear the colored "Actual Results" highlight fills in
#    column G (rows 3-14) back to the plain white look used elsewhere.
# ---------------------------------------------------------------------------
$ws2.Range("F3").Copy()
$ws2.Range("G3").PasteSpecial(-4122)
$ws2.Range("F4").Copy()
$ws2.Range("G4").PasteSpecial(-4122)
$ws2.Range("F5").Copy()
$ws2.Range("G5").PasteSpecial(-4122)
$ws2.Range("F6").Copy()
$ws2.Range("G6").PasteSpecial(-4122)
$ws2.Range("F7").Copy()
$ws2.Range("G7").PasteSpecial(-4122)

$ws2.Range("E8").Copy()
$ws2.Range("G8").PasteSpecial(-4122)
$ws2.Range("E9").Copy()
$ws2.Range("G9").PasteSpecial(-4122)
$ws2.Range("F10").Copy()
$ws2.Range("G10").PasteSpecial(-4122)
$ws2.Range("E11").Copy()
$ws2.Range("G11").PasteSpecial(-4122)
$ws2.Range("E12").Copy()
$ws2.Range("G12").PasteSpecial(-4122)
$ws2.Range("F13").Copy()
$ws2.Range("G13").PasteSpecial(-4122)
$ws2.Range("E14").Copy()
$ws2.Range("G14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) New "Actual Results" text for the patch-applied rows / blank the rest.
# ---------------------------------------------------------------------------
$newReason = "Reason: SR Failed to update. Assign Queue does not exist in OEC. Lookup Assign Queue: cor-CParker."
$ws2.Range("G8").Value = $newReason
$ws2.Range("G9").Value = $newReason
$ws2.Range("G11").Value = $newReason

$ws2.Range("G10").ClearContents()
$ws2.Range("G12").ClearContents()
$ws2.Range("G13").ClearContents()
$ws2.Range("G14").ClearContents()

# ---------------------------------------------------------------------------
# 5) Append the release/version tracking block under the table.
# ---------------------------------------------------------------------------
$ws2.Range("D21").Value = "1.0.37-alpha.4"
$ws2.Range("D22").Value = 49
$ws2.Range("E22").Value = "items"
$ws2.Range("D23").Value = 153
$ws2.Range("E23").Value = "s"
$ws2.Range("D24").Formula = "=D22/D23"

# ---------------------------------------------------------------------------
# 6) View state: keep the freeze pane, restore selections/zoom per sheet and
#    make "RPA Test Cases" the active tab again.
# ---------------------------------------------------------------------------
$ws2.Range("F10").Select()

$ws1.Select()
$ws1.Range("F21").Select()
